$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 460.35715
$ws.Range("J28").Value = 694.6667
$ws.Range("L28").Value = 694.6667
$ws.Range("N28").Value = -1664.6667
$ws.Range("H32").Value = 11365431
$ws.Range("I32").Value = 18182952
$ws.Range("K32").Value = 18182952
$ws.Range("M32").Value = -18182626
$ws.Range("H113").Value = 33336930
$ws.Range("J113").Value = 5169.75
$ws.Range("L113").Value = 5169.75
$ws.Range("N113").Value = -11677.75
$ws.Range("H125").Value = 932.8
$ws.Range("I125").Value = 925.44446
$ws.Range("J125").Value = 999
$ws.Range("K125").Value = 8329.00014
$ws.Range("L125").Value = 8991
$ws.Range("M125").Value = -5869.00014
$ws.Range("N125").Value = -13911
$ws.Range("H135").Value = 1117.5385
$ws.Range("I135").Value = 411.9091
$ws.Range("J135").Value = 4998.5
$ws.Range("K135").Value = 3707.1819
$ws.Range("L135").Value = 44986.5
$ws.Range("M135").Value = -1172.1819
$ws.Range("N135").Value = -50056.5
$ws.Range("H137").Value = 2678.2144
$ws.Range("I137").Value = 3272.7144
$ws.Range("K137").Value = 9818.143199999999
$ws.Range("M137").Value = -7268.143199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1729913.5
$ws.Range("I32").Value = 790334.4
$ws.Range("K32").Value = 790334.4
$ws.Range("M32").Value = -790047.4
$ws.Range("H61").Value = 2905.8667
$ws.Range("I61").Value = 1948.5
$ws.Range("K61").Value = 1948.5
$ws.Range("M61").Value = -1736.5
$ws.Range("H74").Value = 2089.9473
$ws.Range("I74").Value = 1300
$ws.Range("J74").Value = 2967.6667
$ws.Range("K74").Value = 1300
$ws.Range("L74").Value = 2967.6667
$ws.Range("M74").Value = -426
$ws.Range("N74").Value = -4715.6667
$ws.Range("H77").Value = 2089.9473
$ws.Range("I77").Value = 1300
$ws.Range("J77").Value = 2967.6667
$ws.Range("K77").Value = 6500
$ws.Range("L77").Value = 14838.3335
$ws.Range("M77").Value = -2132
$ws.Range("N77").Value = -23574.3335
$ws.Range("H102").Value = 2356.5557
$ws.Range("I102").Value = 1763.7693
$ws.Range("K102").Value = 1763.7693
$ws.Range("M102").Value = -141.7692999999999
$ws.Range("H132").Value = 3156.7334
$ws.Range("I132").Value = 3170.375
$ws.Range("K132").Value = 9511.125
$ws.Range("M132").Value = -6981.125
$ws.Range("H136").Value = 2905.8667
$ws.Range("I136").Value = 1948.5
$ws.Range("K136").Value = 5845.5
$ws.Range("M136").Value = -3295.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 835
$ws.Range("I22").Value = 835
$ws.Range("K22").Value = 835
$ws.Range("M22").Value = -662
$ws.Range("H99").Value = 2429.3076
$ws.Range("I99").Value = 2031.3334
$ws.Range("K99").Value = 2031.3334
$ws.Range("M99").Value = -533.3334
$ws.Range("H134").Value = 3325.3
$ws.Range("I134").Value = 3101.2
$ws.Range("K134").Value = 9303.599999999999
$ws.Range("M134").Value = -6768.599999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H41").Value = 20707.334
$ws.Range("I41").Value = 20707.334
$ws.Range("K41").Value = 20707.334
$ws.Range("M41").Value = -20279.334
$ws.Range("H50").Value = 59497.5
$ws.Range("J50").Value = 59497.5
$ws.Range("L50").Value = 59497.5
$ws.Range("N50").Value = -60747.5
$ws.Range("H58").Value = 1744.8077
$ws.Range("I58").Value = 1192
$ws.Range("J58").Value = 2498.6365
$ws.Range("K58").Value = 1192
$ws.Range("L58").Value = 2498.6365
$ws.Range("M58").Value = -989
$ws.Range("N58").Value = -2904.6365
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H130").Value = 59999.668
$ws.Range("J130").Value = 59999.668
$ws.Range("L130").Value = 59999.668
$ws.Range("N130").Value = -70039.66800000001
$ws.Range("H136").Value = 1744.8077
$ws.Range("I136").Value = 1192
$ws.Range("J136").Value = 2498.6365
$ws.Range("K136").Value = 3576
$ws.Range("L136").Value = 7495.9095
$ws.Range("M136").Value = -1026
$ws.Range("N136").Value = -12595.9095

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 18914
$ws.Range("J74").Value = 17695.5
$ws.Range("L74").Value = 53086.5
$ws.Range("N74").Value = -55208.5
$ws.Range("H77").Value = 18914
$ws.Range("J77").Value = 17695.5
$ws.Range("L77").Value = 159259.5
$ws.Range("N77").Value = -169867.5
$ws.Range("H107").Value = 502.5
$ws.Range("I107").Value = 377
$ws.Range("K107").Value = 1131
$ws.Range("M107").Value = 789
$ws.Range("H108").Value = 12900
$ws.Range("I108").Value = 12900
$ws.Range("K108").Value = 38700
$ws.Range("M108").Value = -35820
$ws.Range("H109").Value = 3474.8572
$ws.Range("I109").Value = 1441.3334
$ws.Range("K109").Value = 4324.0002
$ws.Range("M109").Value = -3284.0002
$ws.Range("H131").Value = 1477211.6
$ws.Range("J131").Value = 1784437.5
$ws.Range("L131").Value = 5353312.5
$ws.Range("N131").Value = -5363392.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1037.7
$ws.Range("J2").Value = 27.833334
$ws.Range("L2").Value = 27.833334
$ws.Range("N2").Value = -253.833334
$ws.Range("H70").Value = 9999.25
$ws.Range("I70").Value = 9998.5
$ws.Range("K70").Value = 9998.5
$ws.Range("M70").Value = -9728.5
$ws.Range("H73").Value = 9999.25
$ws.Range("I73").Value = 9998.5
$ws.Range("K73").Value = 9998.5
$ws.Range("M73").Value = -9062.5
$ws.Range("H122").Value = 2402.3076
$ws.Range("I122").Value = 1634.826
$ws.Range("K122").Value = 4904.478
$ws.Range("M122").Value = -2454.478

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214
$ws.Range("H40").Value = 15608.4
$ws.Range("J40").Value = 4701.3335
$ws.Range("L40").Value = 4701.3335
$ws.Range("N40").Value = -4973.3335
$ws.Range("H61").Value = 1542.5714
$ws.Range("I61").Value = 1466.3334
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1466.3334
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1264.3334
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 1542.5714
$ws.Range("I113").Value = 1466.3334
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1466.3334
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 703.6666
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 4658.3335
$ws.Range("I122").Value = 2103.8572
$ws.Range("J122").Value = 5935.5713
$ws.Range("K122").Value = 6311.571599999999
$ws.Range("L122").Value = 17806.7139
$ws.Range("M122").Value = -3861.571599999999
$ws.Range("N122").Value = -22706.7139
$ws.Range("H136").Value = 7172.75
$ws.Range("I136").Value = 6867.5713
$ws.Range("K136").Value = 20602.7139
$ws.Range("M136").Value = -18052.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15970.5
$ws.Range("I41").Value = 12997
$ws.Range("J41").Value = 18944
$ws.Range("K41").Value = 12997
$ws.Range("L41").Value = 18944
$ws.Range("M41").Value = -12607
$ws.Range("N41").Value = -19724
$ws.Range("H126").Value = 2383.8333
$ws.Range("I126").Value = 2431.3333
$ws.Range("J126").Value = 2336.3333
$ws.Range("K126").Value = 7293.999899999999
$ws.Range("L126").Value = 7008.999899999999
$ws.Range("M126").Value = -4823.999899999999
$ws.Range("N126").Value = -11948.9999
$ws.Range("H136").Value = 3961.6924
$ws.Range("I136").Value = 4045.6365
$ws.Range("K136").Value = 12136.9095
$ws.Range("M136").Value = -9586.9095
